$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 79 (shifts old rows 79/80 down to 80/81)
$ws.Rows.Item(79).Insert()

# Populate the new row 79 with the new weekly data point
$ws.Range("A79").Value2 = 11
$ws.Range("B79").Value2 = "Vega Monumental Concepción"
$ws.Range("C79").Value2 = "Bíobío"
$ws.Range("D79").Value2 = 45267
$ws.Range("D79").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E79").Value2 = 8
$ws.Range("F79").Value2 = 100112026
$ws.Range("G79").Value2 = "Haba"
$ws.Range("H79").Value2 = "Sin especificar"
$ws.Range("I79").Value2 = "Primera"
$ws.Range("J79").Value2 = 50
$ws.Range("K79").Value2 = 12000
$ws.Range("L79").Value2 = 12000
$ws.Range("M79").Value2 = 12000
$ws.Range("N79").Value2 = "$/saco 25 kilos"
$ws.Range("O79").Value2 = "Región del Maule"
$ws.Range("P79").Value2 = 480
$ws.Range("Q79").Value2 = 25
$ws.Range("R79").Value2 = "Hortaliza"
